# Add a "Longitude" column (G) to both station sheets, and fill in the
# still-missing Latitude/Longitude (columns E/F) for the NAF sheet rows
# that didn't have them yet (Bergen, Namsos, Arendal, Steinkjær, Kristiansand).

$wb = $excel.ActiveWorkbook

function Set-TextValue($ws, $addr, $text) {
    # Force the literal text to be stored as a shared string (t="s") rather
    # than being auto-coerced to a number, while leaving the cell's applied
    # style/format untouched (matches cells typed as plain text in Excel).
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

# --- "Viking" sheet: add header G1 = "Longitude" (duplicate of E1) -------
$wsViking = $wb.Worksheets.Item("Viking")
$wsViking.Range("G1").Value = "Longitude"
$wsViking.Range("G1").Font.Bold = $true
$wsViking.Range("G1").Select()

# --- "NAF" sheet: add header G1 = "Longitude" (duplicate of E1) ----------
$wsNaf = $wb.Worksheets.Item("NAF")
$wsNaf.Range("G1").Value = "Longitude"
$wsNaf.Range("G1").Font.Bold = $true

# Fill in the missing Latitude (E) / Longitude (F) values for rows 4-8
Set-TextValue $wsNaf "E4" "60.29552"
Set-TextValue $wsNaf "F4" "5.25276"

Set-TextValue $wsNaf "E5" "64.46625"
Set-TextValue $wsNaf "F5" "11.50387"

Set-TextValue $wsNaf "E6" "58.45858"
Set-TextValue $wsNaf "F6" "8.71815"

Set-TextValue $wsNaf "E7" "64.00874"
Set-TextValue $wsNaf "F7" "11.49363"

# Note: F8 is written before E8 so new shared-string indices land in the
# same order as the authored workbook.
Set-TextValue $wsNaf "F8" "8.12099"
Set-TextValue $wsNaf "E8" "58.17896"

$wsNaf.Range("G1").Select()
